$wb = $excel.ActiveWorkbook

# ALC row 106
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 6179.4546
$ws.Range("I106").Value = 5997.222
$ws.Range("K106").Value = 5997.222
$ws.Range("M106").Value = -5366.222

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2838.125
$ws.Range("I132").Value = 2932.1482
$ws.Range("J132").Value = 2330.4
$ws.Range("K132").Value = 8796.444600000001
$ws.Range("L132").Value = 6991.200000000001
$ws.Range("M132").Value = -6266.444600000001
$ws.Range("N132").Value = -12051.2

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1601.0476
$ws.Range("I137").Value = 992.8461
$ws.Range("K137").Value = 2978.5383
$ws.Range("M137").Value = -428.5383000000002

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3460.5789
$ws.Range("I138").Value = 1301.4706
$ws.Range("K138").Value = 3904.4118
$ws.Range("M138").Value = 1235.5882

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3356.5557
$ws.Range("I61").Value = 3258.1
$ws.Range("J61").Value = 3848.8333
$ws.Range("K61").Value = 3258.1
$ws.Range("L61").Value = 3848.8333
$ws.Range("M61").Value = -3046.1
$ws.Range("N61").Value = -4272.8333

# ARM row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2047
$ws.Range("I88").Value = 1649.5
$ws.Range("K88").Value = 1649.5
$ws.Range("M88").Value = -1243.5

# ARM row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 2047
$ws.Range("I91").Value = 1649.5
$ws.Range("K91").Value = 1649.5
$ws.Range("M91").Value = -245.5

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1833.25
$ws.Range("I110").Value = 1861
$ws.Range("J110").Value = 1750
$ws.Range("K110").Value = 1861
$ws.Range("L110").Value = 1750
$ws.Range("M110").Value = 184
$ws.Range("N110").Value = -5840

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2557.8667
$ws.Range("I132").Value = 2643.923
$ws.Range("K132").Value = 7931.768999999999
$ws.Range("M132").Value = -5401.768999999999

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3356.5557
$ws.Range("I136").Value = 3258.1
$ws.Range("J136").Value = 3848.8333
$ws.Range("K136").Value = 9774.299999999999
$ws.Range("L136").Value = 11546.4999
$ws.Range("M136").Value = -7224.299999999999
$ws.Range("N136").Value = -16646.4999

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2101.4546
$ws.Range("I107").Value = 1651.4546
$ws.Range("K107").Value = 1651.4546
$ws.Range("M107").Value = 268.5454

# CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 725.8570999999999
$ws.Range("I22").Value = 698.5
$ws.Range("K22").Value = 698.5
$ws.Range("M22").Value = -348.5

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5918.5127
$ws.Range("I31").Value = 7661.5654
$ws.Range("K31").Value = 7661.5654
$ws.Range("M31").Value = -7366.5654

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5918.5127
$ws.Range("I34").Value = 7661.5654
$ws.Range("K34").Value = 7661.5654
$ws.Range("M34").Value = -7459.5654

# CRP row 51
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 37498
$ws.Range("J51").Value = 37498
$ws.Range("L51").Value = 37498
$ws.Range("N51").Value = -38970

# CRP row 61
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 37498
$ws.Range("J61").Value = 37498
$ws.Range("L61").Value = 37498
$ws.Range("N61").Value = -38194

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 4577.4287
$ws.Range("I107").Value = 5306.4
$ws.Range("J107").Value = 2755
$ws.Range("K107").Value = 5306.4
$ws.Range("L107").Value = 2755
$ws.Range("M107").Value = -3386.4
$ws.Range("N107").Value = -6595

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4947.5
$ws.Range("I132").Value = 4806.409
$ws.Range("K132").Value = 14419.227
$ws.Range("M132").Value = -11889.227

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3536.8
$ws.Range("I134").Value = 3498.4524
$ws.Range("J134").Value = 3738.125
$ws.Range("K134").Value = 10495.3572
$ws.Range("L134").Value = 11214.375
$ws.Range("M134").Value = -7960.3572
$ws.Range("N134").Value = -16284.375

# CUL row 2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 896.73914
$ws.Range("J2").Value = 372.84616
$ws.Range("L2").Value = 2237.07696
$ws.Range("N2").Value = -2463.07696

# CUL row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 153899340
$ws.Range("I4").Value = 153899340
$ws.Range("K4").Value = 461698020
$ws.Range("M4").Value = -461697908

# CUL row 75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 1521.75
$ws.Range("I75").Value = 2097.5
$ws.Range("J75").Value = 1439.5
$ws.Range("K75").Value = 6292.5
$ws.Range("L75").Value = 4318.5
$ws.Range("M75").Value = -5294.5
$ws.Range("N75").Value = -6314.5

# CUL row 78
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 1521.75
$ws.Range("I78").Value = 2097.5
$ws.Range("J78").Value = 1439.5
$ws.Range("K78").Value = 18877.5
$ws.Range("L78").Value = 12955.5
$ws.Range("M78").Value = -13885.5
$ws.Range("N78").Value = -22939.5

# LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 897.1667
$ws.Range("I16").Value = 897.1667
$ws.Range("K16").Value = 897.1667
$ws.Range("M16").Value = -727.1667

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3342.2307
$ws.Range("I68").Value = 1931
$ws.Range("J68").Value = 5600.2
$ws.Range("K68").Value = 1931
$ws.Range("L68").Value = 5600.2
$ws.Range("M68").Value = -1182
$ws.Range("N68").Value = -7098.2

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 3342.2307
$ws.Range("I71").Value = 1931
$ws.Range("J71").Value = 5600.2
$ws.Range("K71").Value = 9655
$ws.Range("L71").Value = 28001
$ws.Range("M71").Value = -5911
$ws.Range("N71").Value = -35489

# LTW row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1882.56
$ws.Range("I82").Value = 2456.1538
$ws.Range("J82").Value = 1261.1666
$ws.Range("K82").Value = 2456.1538
$ws.Range("L82").Value = 1261.1666
$ws.Range("M82").Value = -2095.1538
$ws.Range("N82").Value = -1983.1666

# LTW row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1882.56
$ws.Range("I85").Value = 2456.1538
$ws.Range("J85").Value = 1261.1666
$ws.Range("K85").Value = 2456.1538
$ws.Range("L85").Value = 1261.1666
$ws.Range("M85").Value = -1208.1538
$ws.Range("N85").Value = -3757.1666

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2371.5264
$ws.Range("I122").Value = 2362.9666
$ws.Range("J122").Value = 2403.625
$ws.Range("K122").Value = 7088.899800000001
$ws.Range("L122").Value = 7210.875
$ws.Range("M122").Value = -4638.899800000001
$ws.Range("N122").Value = -12110.875

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1864.3529
$ws.Range("I136").Value = 1241.2069
$ws.Range("K136").Value = 3723.620699999999
$ws.Range("M136").Value = -1173.620699999999
